$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 193, pushing old rows 193-205 down to 194-206.
$ws.Rows.Item(193).Insert()

# Populate the new row 193 with the new weekly price-report entry.
$ws.Range("A193").Value = 2
$ws.Range("B193").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C193").Value = "Coquimbo"
$ws.Range("D193").Value = 44714
$ws.Range("E193").Value = 4
$ws.Range("F193").Value = 100112031
$ws.Range("G193").Value = "Poroto verde"
$ws.Range("H193").Value = "Magnum"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 500
$ws.Range("K193").Value = 16000
$ws.Range("L193").Value = 18000
$ws.Range("M193").Value = 17000
$ws.Range("N193").Value = "$/malla 25 kilos"
$ws.Range("O193").Value = "Provincia de Limarí"
$ws.Range("P193").Value = 680
$ws.Range("Q193").Value = 25
$ws.Range("R193").Value = "Hortaliza"
